# Webpage Versions.xlsx - add the "Version 2" and "Version 2.1" rows plus
# assorted formatting/view touch-ups, matching the authored update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 4 - Version 2 (2024-09-02), author Kieran Moores
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 2

# Copy the date cell's number formatting (m/d/yyyy) from the row above
# so the new date cell reuses the existing style instead of minting a
# new one, then set its value.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = 45537

$ws.Range("C4").Value = "Kieran Moores"
$ws.Range("D4").Value = "Added logo and improved colour scheme of the website. Added minor functionality and bug fixes."
$ws.Range("E4").Value = "The first iteration of version 2 with a much better, simpler design that is easier to read. Added aesthetics such as a logo and titles"

# ---------------------------------------------------------------------
# Row 5 - Version 2.1 (2024-09-09), author Matthew Adler
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 2.1

$ws.Range("B3").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 45544

$ws.Range("C5").Value = "Matthew Adler"
$ws.Range("D5").Value = "Added a login/logout feature that remembered users. Built a database to handle user accounts and accessibility issues. Started working on a forgot password feature."
$ws.Range("E5").Value = "Improved website functionality that tailors results to the individual. A more personal feel with helpful suggestions."

# ---------------------------------------------------------------------
# Row heights - rows 2 & 3 grow to fit re-wrapped text, new rows 4 & 5
# pick up their own wrapped-text heights.
# ---------------------------------------------------------------------
$ws.Rows(2).RowHeight = 45
$ws.Rows(3).RowHeight = 45
$ws.Rows(4).RowHeight = 60
$ws.Rows(5).RowHeight = 60

# ---------------------------------------------------------------------
# Selection moved to E10
# ---------------------------------------------------------------------
$ws.Range("E10").Select()
